$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update price (D) and volume-1h (E) columns for rows with changes
$ws.Range("D2").Value = "'25.865.64"
$ws.Range("E2").Value = "  -1.49%  "
$ws.Range("D3").Value = "'1.638.61"
$ws.Range("E3").Value = "  -1.19%  "
$ws.Range("E4").Value = "  -0.69%  "
$ws.Range("D5").Value = "'215.64"
$ws.Range("E5").Value = "  -0.49%  "
$ws.Range("D6").Value = "'0.5039"
$ws.Range("E6").Value = "  -1.99%  "
$ws.Range("E7").Value = "  -0.50%  "
$ws.Range("D8").Value = "'0.2570"
$ws.Range("E8").Value = "  -1.31%  "
$ws.Range("D9").Value = "'0.06390"
$ws.Range("E9").Value = "  -1.27%  "
$ws.Range("D10").Value = "'19.70"
$ws.Range("E10").Value = "  -1.39%  "
$ws.Range("D11").Value = "'0.07748"
$ws.Range("E11").Value = "  -1.03%  "
$ws.Range("D12").Value = "'1.650.10"
$ws.Range("E12").Value = "  -0.37%  "
$ws.Range("D13").Value = "'4.258"
$ws.Range("E13").Value = "  -1.30%  "
$ws.Range("D14").Value = "'1.862.82"
$ws.Range("E14").Value = "  -1.25%  "
$ws.Range("D15").Value = "'0.5466"
$ws.Range("E15").Value = "  -1.54%  "
$ws.Range("D16").Value = "'0.0₅7912"
$ws.Range("E16").Value = "  -1.68%  "
$ws.Range("E17").Value = "  -0.52%  "
$ws.Range("D18").Value = "'25.900.23"
$ws.Range("E18").Value = "  -1.39%  "
$ws.Range("D19").Value = "'1.004"
$ws.Range("E19").Value = "  -0.42%  "
$ws.Range("D20").Value = "'202.69"
$ws.Range("E20").Value = "  -3.83%  "
$ws.Range("D21").Value = "'4.392"
$ws.Range("E21").Value = "  -0.79%  "
$ws.Range("D22").Value = "'9.897"
$ws.Range("E22").Value = "  -2.16%  "
$ws.Range("D23").Value = "'5.972"
$ws.Range("E23").Value = "  -1.18%  "
$ws.Range("D24").Value = "'1.005"
$ws.Range("E24").Value = "  -0.35%  "
$ws.Range("D25").Value = "'1.918"
$ws.Range("E25").Value = "  +8.80%  "
$ws.Range("D26").Value = "'140.86"
$ws.Range("E26").Value = "  -2.85%  "
$ws.Range("D27").Value = "'0.1135"
$ws.Range("E27").Value = "  -3.53%  "
$ws.Range("D28").Value = "'15.66"
$ws.Range("E28").Value = "  -1.38%  "
$ws.Range("D29").Value = "'6.760"
$ws.Range("E29").Value = "  -3.57%  "
$ws.Range("D30").Value = "'1.245"
$ws.Range("E30").Value = "  -0.27%  "
$ws.Range("D31").Value = "'0.04974"
$ws.Range("E31").Value = "  -2.86%  "
$ws.Range("D32").Value = "'3.273"
$ws.Range("E32").Value = "  -2.81%  "
$ws.Range("D33").Value = "'3.192"
$ws.Range("E33").Value = "  -1.43%  "
$ws.Range("D34").Value = "'1.547"
$ws.Range("E34").Value = "  -1.43%  "
$ws.Range("E35").Value = "  +0.74%  "
$ws.Range("D36").Value = "'0.8938"
$ws.Range("E36").Value = "  -3.66%  "
$ws.Range("D37").Value = "'2.625"
$ws.Range("E37").Value = "  -4.32%  "
$ws.Range("D38").Value = "'1.150.87"
$ws.Range("E38").Value = "  -2.31%  "
$ws.Range("D39").Value = "'0.5617"
$ws.Range("E39").Value = "  -2.12%  "
$ws.Range("E40").Value = "  -1.52%  "
$ws.Range("D41").Value = "'1.005"
$ws.Range("E41").Value = "  -0.34%  "
$ws.Range("D42").Value = "'5.671"
$ws.Range("E42").Value = "  -0.89%  "
$ws.Range("D45").Value = "'1.774.79"
$ws.Range("E45").Value = "  -1.29%  "
$ws.Range("D46").Value = "'0.0₈118"
$ws.Range("E46").Value = "  +1.15%  "
$ws.Range("D47").Value = "'0.4529"
$ws.Range("E47").Value = "  -0.71%  "
$ws.Range("D48").Value = "'1.004"
$ws.Range("E48").Value = "  -0.36%  "
$ws.Range("D49").Value = "'54.92"
$ws.Range("E49").Value = "  -1.15%  "
$ws.Range("D50").Value = "'0.05053"
$ws.Range("E50").Value = "  -0.67%  "
$ws.Range("D51").Value = "'1.001"
$ws.Range("E51").Value = "  -0.62%  "

# Rows 43/44: TrustWalletToken and Quant swap positions (content swap, A unchanged)
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").Value = "'99.86"
$ws.Range("E43").Value = "  -0.82%  "
$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").Value = "'0.8070"
$ws.Range("E44").Value = "  -2.36%  "
